$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @(
    "Usinfo@sproutup.com",
    "600-4832info@TechWasteRecycling.com",
    "info@techwasterecycling.com",
    "info@arrowrecovery.com",
    "Info@arrowrecovery.com",
    "contact@baytechrecovery.com",
    "info@bmionline.us",
    "moe@datait.com",
    "Moe@datait.com",
    "info@dgglobal.net",
    "contact@evergreenitmanagement.com",
    "info@greenland-resource.com",
    "sales@ironsystems.com",
    "info@ironsystems.com",
    "info@magnakom.com",
    "support@modernwastesolutions.com",
    "onsiterecycling@myoer.com",
    "info@poweron.com",
    "Andy@sem-recycling.com",
    "info@smartwasteusa.com",
    "goldy@starmicro.net",
    "info@t3rs.com"
)

$row = 5
foreach ($email in $emails) {
    $ws.Cells.Item($row, 1).Value = $email
    $row = $row + 1
}
